$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The figures in this sheet are stored as plain text (not numbers), so we
# force Text formatting before assigning the value - otherwise Excel's
# "smart" input would silently convert the numeric-looking string into a
# real number. Resetting the style back to "Normal" afterwards keeps the
# cell formatting identical to the original (General, default style).

# Row 11 - "Enterprises density (per 1000 people)"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "11.85"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "10.45"
$ws.Range("C11").Style = "Normal"

# Row 12 - "Enterprises (% of total)"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "46.13"
$ws.Range("C12").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.43"
$ws.Range("D12").Style = "Normal"
